$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

$rowsData = @(
    @("05", "Ver Usuarios(Bibliotecario)", "", "Ver Apenas Usuarios"),
    @("06", "Ver Usuarios(Admin)", "", "Ver Todos os Usuarios do Sistema"),
    @("07", "Filtrar por Perfil", "Perfil:Usuario", "Ver Apenas Usuarios com Perfil Usuario")
)

foreach ($rowData in $rowsData) {
    $newRow = $table.Rows.Add()
    $newRow.HeadingFormat = 0

    $newRow.Cells.Item(1).Range.Text = $rowData[0]
    $newRow.Cells.Item(2).Range.Text = $rowData[1]
    if ($rowData[2] -ne "") {
        $newRow.Cells.Item(3).Range.Text = $rowData[2]
    }
    $newRow.Cells.Item(4).Range.Text = $rowData[3]
}
